$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new feedback row at row 7
$ws.Range("A7").Value = "rxxx"
$ws.Range("B7").Value = "asdf"
$ws.Range("C7").Value = "asdf"
$ws.Range("D7").Value = "2025-09-30 13:20:59"
